$d = $word.ActiveDocument

# Locate the paragraph that ends with "... Mail, evenementkleur wijzigt, ..."
# (the last paragraph in the document, right before the section break) and
# append two new paragraphs after it, as noted during the meeting.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$newPara1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara1.Range.Text = "- Navragen: Moet een deelnemer bij aanwezigheid inchecken via de mobiele app om te bewijzen dat hij er is?"

$r2 = $newPara1.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()

$newPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara2.Range.Text = "- Moet de administrator de mogelijkheid hebben om een " + [char]0x201C + "anonieme" + [char]0x201D + " user toe te voegen aan een evenement."
